# Add logging (new cache-stat rows/columns) for the "exact method same depth"
# visibility-check test run, mirroring the existing A/I/N/U/AA run columns
# on Sheet2 with a brand-new "AF" run column, and flip cacheEE/cacheEEE to
# true for that run. Also updates the active-sheet/selection view state to
# match where the author left off (Sheet1 selected at J25, Sheet2 scrolled
# over to show the new AF column with AF33 as the last touched cell).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: new "AF" results column -------------------------------------

# Parameter block (rows 1-15) — same settings as the other runs.
$ws2.Range("AF1").Value  = "alldir = false"
$ws2.Range("AF2").Value  = 'dir = "X"'
$ws2.Range("AF3").Value  = "sgn = 1"
$ws2.Range("AF4").Value  = "depth = 6"
$ws2.Range("AF5").Value  = "noSamples = 1000"
$ws2.Range("AF6").Value  = "createRatio = 10"
$ws2.Range("AF7").Value  = "w = 400"
$ws2.Range("AF8").Value  = "h = 400"
$ws2.Range("AF9").Value  = "constructOption = 3"
$ws2.Range("AF10").Value = 'filename = "../models/EEEE_more.obj"'
$ws2.Range("AF11").Value = "width = 800"
$ws2.Range("AF12").Value = "height = 800"
$ws2.Range("AF13").Value = "trace = false"
$ws2.Range("AF14").Value = "sampling = true"
$ws2.Range("AF15").Value = "exact = true"
$ws2.Range("AF18").Value = "cacheCombi = true"

# Results block (rows 20-23) — new cache stats for this run.
$ws2.Range("AF20").Value = "Combi Cache, Size: 837158 Hits: 4642622 Hash Hit Size: 0"
$ws2.Range("AF21").Value = "Edge Edge Edge Cache, Size: 3101 Hits: 53860"
$ws2.Range("AF22").Value = "Edge Edge Cache, Size: 1196 Hits: 21812"
$ws2.Range("AF23").Value = "Completed RST1 in 16481 ms"

# This run is the first to turn both edge-edge and edge-edge-edge caching on.
$ws2.Range("AF16").Value = "cacheEE = true"
$ws2.Range("AF17").Value = "cacheEEE = true"

# --- View state ------------------------------------------------------------

# Sheet2 becomes the background sheet, scrolled right so column K is first
# visible and the new AF column is in view; last selection is AF33.
$ws2.Activate()
$ws2.Range("AF33").Select()
$win = $excel.Application.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1

# Sheet1 is the active tab again, with the last selection at J25.
$ws1.Activate()
$ws1.Range("J25").Select()
